$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.98%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.47%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.304"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.39%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07463"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'10.10%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.808"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.96%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.809"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'10.63%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.458"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.53%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.01%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.01733"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2,578.03%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1696"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'6.58%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07679"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'10.75%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08077"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'6.07%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03008"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.01%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09889"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'10.00%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001496"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.66%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04573"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.88%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006318"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-5.23%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.62%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.12%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.3334"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.59%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'1.93%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.481"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'12.21%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'4.22%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001217"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.04%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004415"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'1.23%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'20.44%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'8.00%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04532"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.85%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007204"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'6.11%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1339"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.00%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002241"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.56%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.12%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006145"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.13%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.873"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.98%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01300"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = "Normal"
